# Updates pulled by the scheduled Sheets runner: refreshed currentAveragePrice /
# LevePrice / LeveProfit figures (columns H-N) for the affected Leve rows across
# the job sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 46687.45
$ws.Range("J17").Value = 49081.527
$ws.Range("L17").Value = 147244.581
$ws.Range("N17").Value = -147580.581
# Row 97: Materia Worth
$ws.Range("H97").Value = 2296.75
$ws.Range("J97").Value = 2296.75
$ws.Range("L97").Value = 6890.25
$ws.Range("N97").Value = -7882.25
# Row 100: Asking for a Friend
$ws.Range("H100").Value = 3527.9412
$ws.Range("I100").Value = 3098.3076
$ws.Range("J100").Value = 4924.25
$ws.Range("K100").Value = 3098.3076
$ws.Range("L100").Value = 4924.25
$ws.Range("M100").Value = -2557.3076
$ws.Range("N100").Value = -6006.25
# Row 113: Amaro Kart
$ws.Range("H113").Value = 5428.857
$ws.Range("J113").Value = 5222.222
$ws.Range("L113").Value = 5222.222
$ws.Range("N113").Value = -11730.222
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1196.9636
$ws.Range("I132").Value = 800.7646999999999
$ws.Range("K132").Value = 2402.2941
$ws.Range("M132").Value = 127.7058999999999
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2259.675
$ws.Range("I138").Value = 1232.55
$ws.Range("J138").Value = 3286.8
$ws.Range("K138").Value = 3697.65
$ws.Range("L138").Value = 9860.400000000001
$ws.Range("M138").Value = 1442.35
$ws.Range("N138").Value = -20140.4
$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 677398.9399999999
$ws.Range("I45").Value = 1266460.1
$ws.Range("K45").Value = 1266460.1
$ws.Range("M45").Value = -1266083.1
# Row 52: Distill and Know that I'm Right
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("N52").Value = 0
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 1484.8654
$ws.Range("I74").Value = 1372.0435
$ws.Range("K74").Value = 1372.0435
$ws.Range("M74").Value = -498.0435
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1484.8654
$ws.Range("I77").Value = 1372.0435
$ws.Range("K77").Value = 6860.2175
$ws.Range("M77").Value = -2492.2175
# Row 102: Smells of Rich Tama-hagane
$ws.Range("H102").Value = 2164.2856
$ws.Range("I102").Value = 2025
$ws.Range("K102").Value = 2025
$ws.Range("M102").Value = -403
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2237.85
$ws.Range("I122").Value = 2019.3125
$ws.Range("J122").Value = 3112
$ws.Range("K122").Value = 6057.9375
$ws.Range("L122").Value = 9336
$ws.Range("M122").Value = -3607.9375
$ws.Range("N122").Value = -14236
$ws = $wb.Worksheets.Item("BSM")
# Row 16: Port of Call: Ul'dah
$ws.Range("H16").Value = 504
$ws.Range("I16").Value = 504
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 504
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = -334
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 9778.409
$ws.Range("I20").Value = 8558.5
$ws.Range("J20").Value = 13031.5
$ws.Range("K20").Value = 8558.5
$ws.Range("L20").Value = 13031.5
$ws.Range("M20").Value = -8311.5
$ws.Range("N20").Value = -13525.5
# Row 47: Lending a Hand
$ws.Range("H47").Value = 154999.75
$ws.Range("J47").Value = 154999.75
$ws.Range("L47").Value = 154999.75
$ws.Range("N47").Value = -156039.75
# Row 53: Kitchen Casualties
$ws.Range("H53").Value = 8000
$ws.Range("J53").Value = 8000
$ws.Range("L53").Value = 8000
$ws.Range("N53").Value = -9148
# Row 94: High Steal
$ws.Range("H94").Value = 549.2759
$ws.Range("I94").Value = 563.3929000000001
$ws.Range("K94").Value = 563.3929000000001
$ws.Range("M94").Value = -112.3929000000001
# Row 104: Hammer and Sails
$ws.Range("H104").Value = 30684
$ws.Range("J104").Value = 30684
$ws.Range("L104").Value = 30684
$ws.Range("N104").Value = -37672
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 5145.3335
$ws.Range("I105").Value = 5773.6665
$ws.Range("K105").Value = 5773.6665
$ws.Range("M105").Value = -4026.6665
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 1485.8948
$ws.Range("I134").Value = 1457.3334
$ws.Range("K134").Value = 4372.0002
$ws.Range("M134").Value = -1837.0002
$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 4165.095
$ws.Range("I31").Value = 1886.2858
$ws.Range("J31").Value = 8722.714
$ws.Range("K31").Value = 1886.2858
$ws.Range("L31").Value = 8722.714
$ws.Range("M31").Value = -1591.2858
$ws.Range("N31").Value = -9312.714
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 4165.095
$ws.Range("I34").Value = 1886.2858
$ws.Range("J34").Value = 8722.714
$ws.Range("K34").Value = 1886.2858
$ws.Range("L34").Value = 8722.714
$ws.Range("M34").Value = -1684.2858
$ws.Range("N34").Value = -9126.714
# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 1958.3529
$ws.Range("I132").Value = 1896.4849
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5689.4547
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3159.4547
$ws.Range("N132").Value = -17060
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 1891.1666
$ws.Range("I134").Value = 1828.5807
$ws.Range("J134").Value = 2279.2
$ws.Range("K134").Value = 5485.742099999999
$ws.Range("L134").Value = 6837.599999999999
$ws.Range("M134").Value = -2950.742099999999
$ws.Range("N134").Value = -11907.6
$ws = $wb.Worksheets.Item("GSM")
# Row 24: Bad Guys Eat Brass
$ws.Range("H24").Value = 19249.5
$ws.Range("I24").Value = 20000
$ws.Range("J24").Value = 18999.334
$ws.Range("K24").Value = 20000
$ws.Range("L24").Value = 18999.334
$ws.Range("M24").Value = -19827
$ws.Range("N24").Value = -19345.334
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 4471.143
$ws.Range("I70").Value = 4262
$ws.Range("K70").Value = 4262
$ws.Range("M70").Value = -3992
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 4471.143
$ws.Range("I73").Value = 4262
$ws.Range("K73").Value = 4262
$ws.Range("M73").Value = -3326
# Row 132: On Board for Lar
$ws.Range("H132").Value = 4004.6
$ws.Range("I132").Value = 3755.75
$ws.Range("K132").Value = 11267.25
$ws.Range("M132").Value = -8737.25
$ws = $wb.Worksheets.Item("LTW")
# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 7839.7896
$ws.Range("I100").Value = 8080.4
$ws.Range("K100").Value = 8080.4
$ws.Range("M100").Value = -7539.4
# Row 122: Hell on Leather
$ws.Range("H122").Value = 8974.777
$ws.Range("I122").Value = 13563.8
$ws.Range("J122").Value = 3238.5
$ws.Range("K122").Value = 40691.39999999999
$ws.Range("L122").Value = 9715.5
$ws.Range("M122").Value = -38241.39999999999
$ws.Range("N122").Value = -14615.5
$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 679.4706
$ws.Range("I136").Value = 596.9375
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = 759.1875
$ws.Range("N136").Value = -11100

# Cells removed entirely for these rows (no longer applicable)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L52").ClearContents()
$ws.Range("M52").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M16").ClearContents()
